$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.015.00'
$ws.Range('E2').Value = '  +3.91%  '
$ws.Range('D3').Value = '2.594.93'
$ws.Range('E3').Value = '  +2.24%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '521.78'
$c.ClearFormats()
$ws.Range('E5').Value = '  +1.71%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '141.07'
$c.ClearFormats()
$ws.Range('E6').Value = '  +0.66%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.996'
$c.ClearFormats()
$ws.Range('E7').Value = '  -0.30%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.567'
$c.ClearFormats()
$ws.Range('E8').Value = '  +2.22%  '
$ws.Range('D9').Value = '2.621.36'
$ws.Range('E9').Value = '  +3.03%  '
$ws.Range('E10').Value = '  +0.21%  '
$ws.Range('E11').Value = '  +1.92%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.332'
$c.ClearFormats()
$ws.Range('E12').Value = '  +2.94%  '
$ws.Range('E13').Value = '  +2.55%  '
$ws.Range('D14').Value = '3.061.54'
$ws.Range('E14').Value = '  +2.57%  '
$ws.Range('D15').Value = '58.982.20'
$ws.Range('E15').Value = '  +3.78%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '20.58'
$c.ClearFormats()
$ws.Range('E16').Value = '  +2.87%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.0000133'
$c.ClearFormats()
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.582.49'
$ws.Range('E18').Value = '  +2.10%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '339.65'
$c.ClearFormats()
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '4.34'
$c.ClearFormats()
$ws.Range('E20').Value = '  +1.59%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '10.21'
$c.ClearFormats()
$ws.Range('E21').Value = '  +1.40%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '6.53'
$c.ClearFormats()
$ws.Range('E22').Value = '  +6.49%  '
$ws.Range('E23').Value = '  -0.27%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '66.26'
$c.ClearFormats()
$ws.Range('E24').Value = '  +2.89%  '
$ws.Range('E25').Value = '  +1.36%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.406'
$c.ClearFormats()
$ws.Range('E26').Value = '  +1.74%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '0.996'
$c.ClearFormats()
$ws.Range('E27').Value = '  -0.34%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '7.13'
$c.ClearFormats()
$ws.Range('E28').Value = '  +3.98%  '
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').Value = '0.0₃0728'
$ws.Range('E30').Value = '  -2.43%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '5.95'
$c.ClearFormats()
$ws.Range('E31').Value = '  -4.65%  '
$ws.Range('E32').Value = '  +1.70%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '18.82'
$c.ClearFormats()
$ws.Range('E33').Value = '  +1.90%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '148.86'
$c.ClearFormats()
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '4.02'
$c.ClearFormats()
$ws.Range('E35').Value = '  +1.16%  '
$ws.Range('E36').Value = '  +0.59%  '
$ws.Range('E37').Value = '  +2.28%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.839'
$c.ClearFormats()
$ws.Range('E38').Value = '  +2.58%  '
$ws.Range('E39').Value = '  +2.75%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.829'
$c.ClearFormats()
$ws.Range('E40').Value = '  -1.25%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '3.56'
$c.ClearFormats()
$ws.Range('E41').Value = '  +2.60%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '277.31'
$c.ClearFormats()
$ws.Range('E42').Value = '  +6.55%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.996'
$c.ClearFormats()
$ws.Range('E43').Value = '  -0.47%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '10.73'
$c.ClearFormats()
$ws.Range('E44').Value = '  +1.08%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.592'
$c.ClearFormats()
$ws.Range('E45').Value = '  +2.54%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.0955'
$c.ClearFormats()
$ws.Range('E46').Value = '  +0.39%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.0521'
$c.ClearFormats()
$ws.Range('E47').Value = '  +0.86%  '
$ws.Range('E48').Value = '  +1.10%  '
$ws.Range('D49').Value = '1.990.09'
$ws.Range('E49').Value = '  +1.28%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '4.62'
$c.ClearFormats()
$ws.Range('E50').Value = '  +2.68%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.0221'
$c.ClearFormats()
$ws.Range('E51').Value = '  +0.33%  '
